$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(81, 1).Value = "Player"
$ws.Cells.Item(81, 2).Value = "د"
$ws.Cells.Item(81, 3).Value = $false
$ws.Cells.Item(81, 4).Value = 0

$ws.Cells.Item(82, 1).Value = "Abdulkarim Almalki"
$ws.Cells.Item(82, 2).Value = "م"
$ws.Cells.Item(82, 3).Value = $false
$ws.Cells.Item(82, 4).Value = 0

$ws.Cells.Item(83, 1).Value = "Abdulkarim Almalki"
$ws.Cells.Item(83, 2).Value = "ا"
$ws.Cells.Item(83, 3).Value = $false
$ws.Cells.Item(83, 4).Value = 0

$ws.Cells.Item(84, 1).Value = "qqqqqqq"
$ws.Cells.Item(84, 2).Value = "م"
$ws.Cells.Item(84, 3).Value = $true
$ws.Cells.Item(84, 4).Value = 2.86

$ws.Cells.Item(85, 1).Value = "Player"
$ws.Cells.Item(85, 2).Value = "ت"
$ws.Cells.Item(85, 3).Value = $true
$ws.Cells.Item(85, 4).Value = 4.06

$ws.Cells.Item(86, 1).Value = "aaaaaa"
$ws.Cells.Item(86, 2).Value = "ل"
$ws.Cells.Item(86, 3).Value = $false
$ws.Cells.Item(86, 4).Value = 0

$ws.Cells.Item(87, 1).Value = "aaaaaa"
$ws.Cells.Item(87, 2).Value = "ج"
$ws.Cells.Item(87, 3).Value = $true
$ws.Cells.Item(87, 4).Value = 2.74

$ws.Cells.Item(88, 1).Value = "aaaaaa"
$ws.Cells.Item(88, 2).Value = "ز"
$ws.Cells.Item(88, 3).Value = $false
$ws.Cells.Item(88, 4).Value = 0

$ws.Cells.Item(89, 1).Value = "EEE"
$ws.Cells.Item(89, 2).Value = "ن"
$ws.Cells.Item(89, 3).Value = $false
$ws.Cells.Item(89, 4).Value = 0

